# Rename column headers: "..._old" -> "..._FV2210", "..._new" -> "..._FV2304"
# then turn the data range into a native Excel Table (ListObject) and freeze
# the header row, matching the target diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($col = 1; $col -le 21; $col++) {
    $cell = $ws.Cells.Item(1, $col)
    $text = [string]$cell.Value2
    if ($text -like "*_old") {
        $cell.Value2 = ($text -replace "_old$", "_FV2210")
    } elseif ($text -like "*_new") {
        $cell.Value2 = ($text -replace "_new$", "_FV2304")
    }
}

# Convert the used range A1:U56 into a table with an autofilter.
$range = $ws.Range("A1:U56")
$table = $ws.ListObjects.Add(1, $range, $null, 1)
$table.Name = "Table1"

# Freeze the header row (row 1).
$ws.Activate()
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
